$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.049.85"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.44%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.382.15"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "557.72"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.23%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.12"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.99%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("E8").Value = "  -0.63%  "

$ws.Range("E9").Value = "  +1.23%  "

$ws.Range("E10").Value = "  -0.11%  "

$ws.Range("E12").Value = "  -2.73%  "

$ws.Range("E13").Value = "  -3.13%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.809.27"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.93%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "60.000.63"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.46%  "

$ws.Range("E16").Value = "  +0.87%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.385.90"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.56%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.13"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.61%  "

$ws.Range("E19").Value = "  +2.07%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "321.90"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.97%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.67"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.29%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("D22").Style = "Normal"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "64.13"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.53%  "

$ws.Range("E24").Value = "  +0.58%  "

$ws.Range("E25").Value = "  -0.08%  "

$ws.Range("E26").Value = "  -2.36%  "

$ws.Range("E27").Value = "  +0.96%  "

$ws.Range("E28").Value = "  +2.50%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0763"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.94%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "170.11"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.96%  "

$ws.Range("E31").Value = "  +1.34%  "

$ws.Range("E32").Value = "  +10.61%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.400"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.19%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "18.17"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.28%  "

$ws.Range("E35").Value = "  +2.09%  "

$ws.Range("E37").Value = "  +0.08%  "

$ws.Range("E38").Value = "  -0.60%  "

$ws.Range("E39").Value = "  -0.38%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "320.25"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.38%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "38.65"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.75%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "145.34"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.24%  "

$ws.Range("E43").Value = "  -2.97%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0972"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.52%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "19.84"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.61%  "

$ws.Range("E46").Value = "  -0.13%  "

$ws.Range("E47").Value = "  -0.99%  "

$ws.Range("E48").Value = "  -1.88%  "

$ws.Range("E49").Value = "  +0.26%  "

$ws.Range("E50").Value = "  -0.95%  "

$ws.Range("E51").Value = "  -0.06%  "
